# Update cryptos list with latest price and volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.378.90"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.686.41"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.53"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5532"
$ws.Range("E6").Value = "  +8.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.009"
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06504"
$ws.Range("E9").Value = "  +1.77%  "
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07577"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.553"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.685.20"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5819"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008471"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.46"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.418.81"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.943"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.37"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.009"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.75"
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1325"
$ws.Range("E25").Value = "  +10.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.908"
$ws.Range("E26").Value = "  +3.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.83"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06323"
$ws.Range("E28").Value = "  -4.57%  "
$ws.Range("E29").Value = "  +4.39%  "
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.599"
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.673"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("E34").Value = "  +2.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6235"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.408"
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.238"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.116.92"
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01634"
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8774"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.65"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.836.13"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("E45").Value = "  -4.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.42"
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.209"
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05280"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4301"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.087"
$ws.Range("E51").Value = "  +0.78%  "
